$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chapter1-2")
$ws.Activate()

# Update dialogue text cells. Order matters: new/unique strings must be
# assigned in the same sequence as they first appear, so the rebuilt shared
# string table lands in the same order as the target workbook.
$ws.Range("C12").Value = '聽起來蠻酷炫的，實際上是怎麼運作的呢？'
$ws.Range("C8").Value = '你想要從那邊開始了解起呢？'
$ws.Range("C11").Value = '在套用這套系統框架之後，遊戲開發者就能夠隨意的在遊戲中加入各種對話演出與劇情。'
$ws.Range("C22").Value = '例如：文字或圖片的移動、旋轉，或者是淡出跟淡入。'
$ws.Range("C23").Value = '好像有點複雜，有什麼快速入門的方法呢？'
$ws.Range("C26").Value = '真令人期待！'
$ws.Range("C5").Value = '是程式大大！正好有個問題想請教你。'
$ws.Range("C2").Value = '第一章： Horizon Dialogue Plugin 簡介'
$ws.Range("C6").Value = '是想問這個 Horizon Dialogue Plugin 該怎麼使用嗎？'
$ws.Range("C7").Value = '真不愧是大大！馬上就猜出我想問什麼了！剛剛才拿到這個 Plugin ，正困擾著該怎麼使用呢。'
$ws.Range("C9").Value = '可以先談談 Plugin 主要的用途是什麼嗎？'
$ws.Range("C10").Value = '簡單來說，這個 Plugin 的主要目標是想要建立一套遊戲系統框架。'
$ws.Range("C13").Value = '首先，在開始使用這套 Plugin 之前，我建議先去研究 HorizonUI 、 HorizonTween 以及 HorizonFramework 這幾個 Plugin 的功能該怎麼使用。'
$ws.Range("C14").Value = '因為 Dialogue Plugin 非常緊密的整合了這幾套 Plugin ，並在增加了對話演出所需要的流程與功能。'
$ws.Range("C15").Value = '剛剛所提到的遊戲系統框架，跟 HorizonFramework 有什麼關係嗎？'
$ws.Range("C16").Value = '是的， HorizonDialogue Plugin 的系統框架就是對 HorizonFramework 進行擴增，二套 Plugin 在工作流程的設計基本上是一致的。'
$ws.Range("C17").Value = '只是為了加入對話事件，因此在 HorizonScene 的生命週期中新插入了： AddDialogueEvent 這個事件。'
$ws.Range("C18").Value = '嗯，大致上有那麼一點感覺了。所以說 HorizonUI 的用途主要是用來顯示對話的囉？'
$ws.Range("C19").Value = '是的，藉由 HorizonUI 中所提供的功能，就能夠實現多樣的文字對話效果。'
$ws.Range("C20").Value = '那 HorizonTween 呢？'
$ws.Range("C21").Value = 'HorizonTween 的話，主要是拿來實現一些文字與圖片的動畫效果。'
$ws.Range("C24").Value = '我想，先到 GitHub 上把這個專案 Demo 抓下來玩玩是最快的方法了。'
$ws.Range("C25").Value = '在這之後我會示範幾個 Plugin 的進階用法。'

$ws.Range("C22").Select()

